$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Coby White"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Chicago Bulls"

$ws.Range("A3").Value = "Trae Young"
$ws.Range("B3").Value = "PG"
$ws.Range("C3").Value = "Atlanta Hawks"

$ws.Range("A4").Value = "P.J. Washington Jr."
$ws.Range("B4").Value = "PF"
$ws.Range("C4").Value = "Dallas Mavericks"

$ws.Range("A5").Value = "Shaedon Sharpe"
$ws.Range("B5").Value = "SG,SF"
$ws.Range("C5").Value = "Portland Trail Blazers"

$ws.Range("A6").Value = "Jabari Smith Jr."
$ws.Range("B6").Value = "PF,C"
$ws.Range("C6").Value = "Houston Rockets"

$ws.Range("A7").Value = "Alperen Sengün"
$ws.Range("B7").Value = "C"
$ws.Range("C7").Value = "Houston Rockets"

$ws.Range("A8").Value = "Dereck Lively II"
$ws.Range("B8").Value = "C"
$ws.Range("C8").Value = "Dallas Mavericks"

$ws.Range("A9").Value = "Nicolas Claxton"
$ws.Range("B9").Value = "C"
$ws.Range("C9").Value = "Brooklyn Nets"

$ws.Range("A10").Value = "Devin Booker"
$ws.Range("B10").Value = "PG,SG"
$ws.Range("C10").Value = "Phoenix Suns"

$ws.Range("A12").Value = "Jalen Brunson"
$ws.Range("B12").Value = "PG"
$ws.Range("C12").Value = "New York Knicks"

$ws.Range("A13").Value = "LeBron James"
$ws.Range("B13").Value = "SF,PF"
$ws.Range("C13").Value = "Los Angeles Lakers"

$ws.Range("A14").Value = "Desmond Bane"
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "Memphis Grizzlies"

$ws.Range("A15").Value = "Walker Kessler"
$ws.Range("B15").Value = "C"
$ws.Range("C15").Value = "Utah Jazz"

$ws.Range("A16").Value = "Devin Vassell"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "San Antonio Spurs"

$ws.Range("A19").Value = "Norman Powell"
$ws.Range("B19").Value = "SG,SF"
$ws.Range("C19").Value = "LA Clippers"
